$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Daily update from Airtable: a new session record is appended as row 14
# (previously an empty spacer row). The row reproduces the same look and
# feel as the other data rows (e.g. row 13), with its own set of values.
# ---------------------------------------------------------------------

# --- Values -------------------------------------------------------------
$ws.Range("A14").Value = "BARCELONA ACTIVA SA SOCIETAT PRIVADA MUNICIPAL"
$ws.Range("B14").Value = "Orientació Barcelona Activa 2024"
$ws.Range("C14").Value = "24/BCNACTIVA24/2024"
$ws.Range("D14").Value = "Tutoria"
$ws.Range("E14").Value = "Pere Girona Campi"
$ws.Range("F14").Value = 45606.584502314814
$ws.Range("G14").Value = 45624
$ws.Range("H14").Value = 45624
$ws.Range("I14").Value = 25569.479166666668
$ws.Range("J14").Value = 25569.5625
$ws.Range("K14").Value = "Individual"
$ws.Range("L14").Value = "Orientació integral generalista"
$ws.Range("M14").Value = "Diagnòstic inicial "
$ws.Range("N14").Value = "Virtual"
$ws.Range("O14").Value = "43453863D - PEPITO MENGANITO JUANITO"
$ws.Range("P14").Value = "43453863D-11-2024-2"
$ws.Range("Q14").Value = "En curs"
$ws.Range("R14").Value = "GARANTIA JUVENIL"
$ws.Range("S14").Value = "NO_APLICA"
$ws.Range("T14").Value = "NO_APLICA"
$ws.Range("U14").Value = "Definitiu"
$ws.Range("V14").Value = ""

# --- Row height (matches the other data rows) ---------------------------
$ws.Rows("14").RowHeight = 63.75

# --- Formatting -----------------------------------------------------------
# Common look shared by every populated cell in the row: Arial 10, not
# bold/italic/underlined, centered, wrapped text.
$dataRange = $ws.Range("A14:V14")
$dataRange.Font.Name = "Arial"
$dataRange.Font.Size = 10
$dataRange.Font.Bold = $false
$dataRange.Font.Italic = $false
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true

# Date columns (F, G, H)
$ws.Range("F14:H14").NumberFormat = "dd/MM/yyyy"

# Time columns (I, J)
$ws.Range("I14:J14").NumberFormat = "HH:mm:ss"

# Column O (Process) is left/top aligned instead of centered.
$ws.Range("O14").HorizontalAlignment = -4131
$ws.Range("O14").VerticalAlignment = -4160

# A14 carries a left border like the first column of every other data row
# (xlEdgeLeft = 7).
$ws.Range("A14").Borders.Item(7).LineStyle = 1
$ws.Range("A14").Borders.Item(7).Weight = 2

# Trailing helper columns (W:AH) stay blank, matching the rest of the table
# (still present as styled-but-empty cells).
$ws.Range("W14:AH14").Locked = $true
